$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.517.77'
$ws.Range('D3').Value = '1.728.12'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.57'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4808'
$ws.Range('E7').Value = '  +1.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2670'
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').Value = '1.735.51'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07197'
$ws.Range('E11').Value = '  +1.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.59'
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6115'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.527'
$ws.Range('E14').Value = '  +2.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.19'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9994'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '26.532.75'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9993'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006942'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.56'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '1.956.17'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.524'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.796'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.248'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.98'
$ws.Range('E25').Value = '  +1.57%  '
$ws.Range('E26').Value = '  +1.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.777'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.397'
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '107.25'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.963'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08022'
$ws.Range('E31').Value = '  +3.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.691'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04523'
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9987'
$ws.Range('E35').Value = '  +2.53%  '
$ws.Range('E36').Value = '  +1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.079'
$ws.Range('E37').Value = '  +8.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9111'
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.371'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.003'
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '103.31'
$ws.Range('E41').Value = '  -9.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01506'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.635'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.984'
$ws.Range('E45').Value = '  +11.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1182'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05357'
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.791'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.48'
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.255'
$ws.Range('E50').Value = '  +3.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.27'
$ws.Range('E51').Value = '  +1.26%  '
